$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1217.6
$arr[0,1] = 298.5
$arr[0,2] = 1359
$arr[0,3] = 895.5
$arr[0,4] = 4077
$arr[0,5] = -727.5
$arr[0,6] = -4413
$ws.Range("H17:N17").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2933.3333
$arr[0,1] = 2933.3333
$arr[0,2] = 0
$arr[0,3] = 2933.3333
$arr[0,4] = 0
$arr[0,5] = -2607.3333
$arr[0,6] = $null
$ws.Range("H32:N32").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 17336
$arr[0,1] = 0
$arr[0,2] = 17336
$arr[0,3] = 0
$arr[0,4] = 52008
$arr[0,5] = $null
$arr[0,6] = -52246
$ws.Range("H46:N46").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 17336
$arr[0,1] = 0
$arr[0,2] = 17336
$arr[0,3] = 0
$arr[0,4] = 52008
$arr[0,5] = $null
$arr[0,6] = -52976
$ws.Range("H60:N60").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 728
$arr[0,1] = 611.5
$arr[0,2] = 883.3333
$arr[0,3] = 1834.5
$arr[0,4] = 2649.9999
$arr[0,5] = 1232.5
$arr[0,6] = -8783.999899999999
$ws.Range("H111:N111").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3036324.8
$arr[0,1] = 5634.857
$arr[0,2] = 8340032
$arr[0,3] = 5634.857
$arr[0,4] = 8340032
$arr[0,5] = -2192.857
$arr[0,6] = -8346916
$ws.Range("H116:N116").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1311.8889
$arr[0,1] = 1101
$arr[0,2] = 2999
$arr[0,3] = 3303
$arr[0,4] = 8997
$arr[0,5] = 1837
$arr[0,6] = -19277
$ws.Range("H138:N138").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 94460
$arr[0,1] = 0
$arr[0,2] = 94460
$arr[0,3] = 0
$arr[0,4] = 94460
$arr[0,5] = $null
$arr[0,6] = -104740
$ws.Range("H139:N139").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6539.868
$arr[0,1] = 2906.535
$arr[0,2] = 22163.2
$arr[0,3] = 2906.535
$arr[0,4] = 22163.2
$arr[0,5] = -2619.535
$arr[0,6] = -22737.2
$ws.Range("H32:N32").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 52775
$arr[0,1] = 0
$arr[0,2] = 52775
$arr[0,3] = 0
$arr[0,4] = 52775
$arr[0,5] = $null
$arr[0,6] = -62595
$ws.Range("H124:N124").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 68000
$arr[0,1] = 0
$arr[0,2] = 68000
$arr[0,3] = 0
$arr[0,4] = 68000
$arr[0,5] = $null
$arr[0,6] = -77960
$ws.Range("H128:N128").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2121.1516
$arr[0,1] = 1928.5358
$arr[0,2] = 3199.8
$arr[0,3] = 5785.607400000001
$arr[0,4] = 9599.400000000001
$arr[0,5] = -3255.607400000001
$arr[0,6] = -14659.4
$ws.Range("H132:N132").Value = $arr

$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1498.3334
$arr[0,1] = 1060.625
$arr[0,2] = 5000
$arr[0,3] = 1060.625
$arr[0,4] = 5000
$arr[0,5] = -576.625
$arr[0,6] = -5968
$ws.Range("H54:N54").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 51041.383
$arr[0,1] = 65370.562
$arr[0,2] = 5188
$arr[0,3] = 65370.562
$arr[0,4] = 5188
$arr[0,5] = -63623.562
$arr[0,6] = -8682
$ws.Range("H105:N105").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 32977.824
$arr[0,1] = 0
$arr[0,2] = 32977.824
$arr[0,3] = 0
$arr[0,4] = 32977.824
$arr[0,5] = $null
$arr[0,6] = -43097.824
$ws.Range("H132:N132").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 8183.6
$arr[0,1] = 3717.6
$arr[0,2] = 12649.6
$arr[0,3] = 11152.8
$arr[0,4] = 37948.8
$arr[0,5] = -8617.799999999999
$arr[0,6] = -43018.8
$ws.Range("H134:N134").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 88890.71000000001
$arr[0,1] = 0
$arr[0,2] = 88890.71000000001
$arr[0,3] = 0
$arr[0,4] = 88890.71000000001
$arr[0,5] = $null
$arr[0,6] = -99170.71000000001
$ws.Range("H138:N138").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 43499.06
$arr[0,1] = 0
$arr[0,2] = 43499.06
$arr[0,3] = 0
$arr[0,4] = 43499.06
$arr[0,5] = $null
$arr[0,6] = -53859.06
$ws.Range("H140:N140").Value = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1500
$arr[0,1] = 1500
$arr[0,2] = 0
$arr[0,3] = 1500
$arr[0,4] = 0
$arr[0,5] = -1326
$arr[0,6] = $null
$ws.Range("H17:N17").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3732.55
$arr[0,1] = 2422.875
$arr[0,2] = 4605.6665
$arr[0,3] = 2422.875
$arr[0,4] = 4605.6665
$arr[0,5] = -2127.875
$arr[0,6] = -5195.6665
$ws.Range("H31:N31").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3732.55
$arr[0,1] = 2422.875
$arr[0,2] = 4605.6665
$arr[0,3] = 2422.875
$arr[0,4] = 4605.6665
$arr[0,5] = -2220.875
$arr[0,6] = -5009.6665
$ws.Range("H34:N34").Value = $arr

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2502124.8
$arr[0,1] = 499
$arr[0,2] = 3336000
$arr[0,3] = 1497
$arr[0,4] = 10008000
$arr[0,5] = -1159
$arr[0,6] = -10008676
$ws.Range("H41:N41").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4430.6665
$arr[0,1] = 4430.6665
$arr[0,2] = 0
$arr[0,3] = 13291.9995
$arr[0,4] = 0
$arr[0,5] = -13021.9995
$arr[0,6] = $null
$ws.Range("H64:N64").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4430.6665
$arr[0,1] = 4430.6665
$arr[0,2] = 0
$arr[0,3] = 13291.9995
$arr[0,4] = 0
$arr[0,5] = -12355.9995
$arr[0,6] = $null
$ws.Range("H67:N67").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 114.09091
$arr[0,1] = 110.125
$arr[0,2] = 124.666664
$arr[0,3] = 330.375
$arr[0,4] = 373.999992
$arr[0,5] = 165.625
$arr[0,6] = -1365.999992
$ws.Range("H97:N97").Value = $arr

$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 174.07143
$arr[0,1] = 38.333332
$arr[0,2] = 211.09091
$arr[0,3] = 38.333332
$arr[0,4] = 211.09091
$arr[0,5] = 74.666668
$arr[0,6] = -437.09091
$ws.Range("H2:N2").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6945.75
$arr[0,1] = 3900
$arr[0,2] = 9991.5
$arr[0,3] = 3900
$arr[0,4] = 9991.5
$arr[0,5] = -3607
$arr[0,6] = -10577.5
$ws.Range("H18:N18").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 8481.286
$arr[0,1] = 14599.5
$arr[0,2] = 6034
$arr[0,3] = 14599.5
$arr[0,4] = 6034
$arr[0,5] = -14340.5
$arr[0,6] = -6552
$ws.Range("H52:N52").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 8373.5
$arr[0,1] = 10122.25
$arr[0,2] = 6624.75
$arr[0,3] = 10122.25
$arr[0,4] = 6624.75
$arr[0,5] = -9852.25
$arr[0,6] = -7164.75
$ws.Range("H70:N70").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 8373.5
$arr[0,1] = 10122.25
$arr[0,2] = 6624.75
$arr[0,3] = 10122.25
$arr[0,4] = 6624.75
$arr[0,5] = -9186.25
$arr[0,6] = -8496.75
$ws.Range("H73:N73").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 466.5
$arr[0,1] = 360
$arr[0,2] = 999
$arr[0,3] = 360
$arr[0,4] = 999
$arr[0,5] = 136
$arr[0,6] = -1991
$ws.Range("H97:N97").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 52799.2
$arr[0,1] = 0
$arr[0,2] = 52799.2
$arr[0,3] = 0
$arr[0,4] = 52799.2
$arr[0,5] = $null
$arr[0,6] = -57699.2
$ws.Range("H123:N123").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6148.0713
$arr[0,1] = 5232.4
$arr[0,2] = 8437.25
$arr[0,3] = 15697.2
$arr[0,4] = 25311.75
$arr[0,5] = -13167.2
$arr[0,6] = -30371.75
$ws.Range("H132:N132").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 60250
$arr[0,1] = 0
$arr[0,2] = 60250
$arr[0,3] = 0
$arr[0,4] = 60250
$arr[0,5] = $null
$arr[0,6] = -70530
$ws.Range("H139:N139").Value = $arr

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5464.75
$arr[0,1] = 4999
$arr[0,2] = 5531.2856
$arr[0,3] = 4999
$arr[0,4] = 5531.2856
$arr[0,5] = -4887
$arr[0,6] = -5755.2856
$ws.Range("H7:N7").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 15460
$arr[0,1] = 0
$arr[0,2] = 15460
$arr[0,3] = 0
$arr[0,4] = 15460
$arr[0,5] = $null
$arr[0,6] = -25136
$ws.Range("H119:N119").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5464.75
$arr[0,1] = 4999
$arr[0,2] = 5531.2856
$arr[0,3] = 14997
$arr[0,4] = 16593.8568
$arr[0,5] = -12527
$arr[0,6] = -21533.8568
$ws.Range("H126:N126").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3101.2778
$arr[0,1] = 3242.5881
$arr[0,2] = 699
$arr[0,3] = 9727.764299999999
$arr[0,4] = 2097
$arr[0,5] = -7197.764299999999
$arr[0,6] = -7157
$ws.Range("H132:N132").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1737.6207
$arr[0,1] = 1495.6957
$arr[0,2] = 2665
$arr[0,3] = 4487.0871
$arr[0,4] = 7995
$arr[0,5] = -1937.0871
$arr[0,6] = -13095
$ws.Range("H136:N136").Value = $arr

$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10272.818
$arr[0,1] = 1125.125
$arr[0,2] = 34666.668
$arr[0,3] = 2250.25
$arr[0,4] = 69333.336
$arr[0,5] = -1189.25
$arr[0,6] = -71455.336
$ws.Range("H81:N81").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10272.818
$arr[0,1] = 1125.125
$arr[0,2] = 34666.668
$arr[0,3] = 11251.25
$arr[0,4] = 346666.68
$arr[0,5] = -5947.25
$arr[0,6] = -357274.68
$ws.Range("H84:N84").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 25000
$arr[0,1] = 0
$arr[0,2] = 25000
$arr[0,3] = 0
$arr[0,4] = 25000
$arr[0,5] = $null
$arr[0,6] = -31988
$ws.Range("H105:N105").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1786.826
$arr[0,1] = 1199.8823
$arr[0,2] = 3449.8333
$arr[0,3] = 3599.6469
$arr[0,4] = 10349.4999
$arr[0,5] = -1679.6469
$arr[0,6] = -14189.4999
$ws.Range("H107:N107").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 38000
$arr[0,1] = 0
$arr[0,2] = 38000
$arr[0,3] = 0
$arr[0,4] = 38000
$arr[0,5] = $null
$arr[0,6] = -41494
$ws.Range("H121:N121").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 60000
$arr[0,1] = 0
$arr[0,2] = 60000
$arr[0,3] = 0
$arr[0,4] = 60000
$arr[0,5] = $null
$arr[0,6] = -69840
$ws.Range("H125:N125").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2264.5
$arr[0,1] = 2125.087
$arr[0,2] = 3333.3333
$arr[0,3] = 6375.261
$arr[0,4] = 9999.999899999999
$arr[0,5] = -3845.261
$arr[0,6] = -15059.9999
$ws.Range("H132:N132").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 139925.67
$arr[0,1] = 0
$arr[0,2] = 139925.67
$arr[0,3] = 0
$arr[0,4] = 139925.67
$arr[0,5] = $null
$arr[0,6] = -150125.67
$ws.Range("H137:N137").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 150000
$arr[0,1] = 0
$arr[0,2] = 150000
$arr[0,3] = 0
$arr[0,4] = 150000
$arr[0,5] = $null
$arr[0,6] = -150000
$ws.Range("H138:N138").Value = $arr

Write-Host "All edits applied."
